# A new tracking entry was logged on 2024-09-07 12:12:22 for the "2024"
# sheet's September column (R/S). The sheet is laid out newest-first, so
# the new entry is inserted as a brand-new row 35, pushing every
# following row (old 35 -> new 36, ..., old 91 -> new 92) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a blank row above row 35; this shifts rows 35:91 down to 36:92
# and extends the sheet dimension from A1:Y91 to A1:Y92.
$ws.Rows(35).Insert()

# Populate the newly inserted row 35 with the new September log entry.
$ws.Range("R35").Value = "balance your axis"
$ws.Range("S35").Value = "2024-09-07 12:12:22"
